# MemorialGander 2018 - XLS per attrezzo finale (final apparatus start order)
# Re-seeds the Men ("M") and Women ("F") final-apparatus starting-order lists
# with the finalised draw, restyles the "M+F" combined-view Women block to
# left-aligned text, and leaves the cursor parked where the author left it on
# each tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "M" - Men final apparatus starting order (rows 5-14: Name / Apparatus)
# ---------------------------------------------------------------------
$wsM = $wb.Worksheets.Item("M")

$menOrder = @(
  @("Marian Dragulescu", "Pommel Horse"),
  @("Pablo Brägger", "Horizontal Bar"),
  @("Julien Gobaux", "Rings"),
  @("Oliver Hegi", "Vault"),
  @("Nikita Nagornyy", "Vault"),
  @("Marcel Nguyen", "Vault"),
  @("Bart Deurloo", "Rings"),
  @("Oleg Verniaiev", "Parallel Bars"),
  @("Arthur Nory Oyakawa Mariano", "Parallel Bars"),
  @("Cory Paterson", "Floor")
)

$r = 5
foreach ($row in $menOrder) {
  $wsM.Cells.Item($r, 2).Value = $row[0]
  $wsM.Cells.Item($r, 3).Value = $row[1]
  $r++
}

$wsM.Activate()
$wsM.Range("A11").Select()

# ---------------------------------------------------------------------
# Sheet "F" - Women final apparatus starting order (rows 5-15: Name / Apparatus)
# ---------------------------------------------------------------------
$wsF = $wb.Worksheets.Item("F")

$womenOrder = @(
  @("Angelina Melnikova", "Vault"),
  @("Kim Bui", "Floor"),
  @("Flávia Lopes Saraiva ", "Vault"),
  @("Eythora Thorsdottir", "Floor"),
  @("Diana Varinska", "Floor"),
  @("Jade Barbosa", "Floor"),
  @("Caterina Barloggio", "Floor"),
  @("Anina Wildi", "Uneven bars"),
  @("Oksana Chusovitina", "Balance Beam"),
  @("Ilaria Kaeslin", "Floor"),
  @("Mélanie De Jesus Dos Santos", "Floor")
)

$r = 5
foreach ($row in $womenOrder) {
  $wsF.Cells.Item($r, 2).Value = $row[0]
  $wsF.Cells.Item($r, 3).Value = $row[1]
  $r++
}

$wsF.Activate()
$wsF.Range("A12").Select()

# ---------------------------------------------------------------------
# Sheet "M+F" - combined overview; Women block (B23:C28) switches from
# centered to left-aligned text. Values recompute from the M!/F! formulas.
# ---------------------------------------------------------------------
$wsMF = $wb.Worksheets.Item("M+F")
$wsMF.Range("B23:C28").HorizontalAlignment = -4131

$wsMF.Activate()
$wsMF.Range("G4").Select()
